$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.714.01"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").Value = "1.645.83"
$ws.Range("E3").Value = "  +0.75%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'214.31"
$ws.Range("E5").Value = "  +0.57%  "
$ws.Range("D6").Value = "'0.505"
$ws.Range("E6").Value = "  +1.99%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  -0.70%  "
$ws.Range("D9").Value = "'0.0628"
$ws.Range("E9").Value = "  +0.34%  "
$ws.Range("D10").Value = "'19.16"
$ws.Range("E10").Value = "  +0.70%  "
$ws.Range("E11").Value = "  -0.10%  "
$ws.Range("D12").Value = "1.865.66"
$ws.Range("E12").Value = "  +0.23%  "
$ws.Range("D13").Value = "1.648.12"
$ws.Range("E13").Value = "  +0.82%  "
$ws.Range("E14").Value = "  +1.89%  "
$ws.Range("D15").Value = "'0.531"
$ws.Range("E15").Value = "  +0.79%  "
$ws.Range("D16").Value = "'65.64"
$ws.Range("E16").Value = "  +3.86%  "
$ws.Range("D17").Value = "26.725.06"
$ws.Range("E17").Value = "  +0.27%  "
$ws.Range("D18").Value = "0.0₃0746"
$ws.Range("E18").Value = "  +0.81%  "
$ws.Range("D19").Value = "'218.22"
$ws.Range("E19").Value = "  +3.64%  "
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("D21").Value = "'4.35"
$ws.Range("E21").Value = "  +0.94%  "
$ws.Range("D22").Value = "'6.33"
$ws.Range("E22").Value = "  +1.75%  "
$ws.Range("D23").Value = "'9.41"
$ws.Range("E23").Value = "  -0.49%  "
$ws.Range("D24").Value = "'2.14"
$ws.Range("E24").Value = "  +11.13%  "
$ws.Range("D25").Value = "'147.54"
$ws.Range("E25").Value = "  +0.20%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").Value = "'0.121"
$ws.Range("E27").Value = "  -0.46%  "
$ws.Range("D28").Value = "'6.95"
$ws.Range("E28").Value = "  +0.87%  "
$ws.Range("D29").Value = "'15.77"
$ws.Range("E29").Value = "  +2.36%  "
$ws.Range("D30").Value = "'0.0519"
$ws.Range("E30").Value = "  -0.74%  "
$ws.Range("D31").Value = "'1.17"
$ws.Range("E31").Value = "  -0.52%  "
$ws.Range("E32").Value = "  +3.93%  "
$ws.Range("E33").Value = "  +2.05%  "
$ws.Range("D34").Value = "1.272.04"
$ws.Range("E34").Value = "  +8.65%  "
$ws.Range("E35").Value = "  +0.87%  "
$ws.Range("E36").Value = "  +0.98%  "
$ws.Range("D38").Value = "'0.810"
$ws.Range("E38").Value = "  -0.22%  "
$ws.Range("D39").Value = "'0.514"
$ws.Range("E39").Value = "  +1.52%  "
$ws.Range("E40").Value = "  -0.06%  "
$ws.Range("E41").Value = "  -1.62%  "
$ws.Range("D42").Value = "'0.803"
$ws.Range("E42").Value = "  +1.14%  "
$ws.Range("E43").Value = "  -0.22%  "
$ws.Range("D44").Value = "1.776.05"
$ws.Range("E44").Value = "  +0.34%  "
$ws.Range("D45").Value = "'93.88"
$ws.Range("E45").Value = "  +1.46%  "
$ws.Range("D46").Value = "'1.61"
$ws.Range("E46").Value = "  +3.83%  "
$ws.Range("D47").Value = "'55.50"
$ws.Range("E47").Value = "  +1.45%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "'0.0513"
$ws.Range("E48").Value = "  +0.36%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'7.62"
$ws.Range("E49").Value = "  +0.83%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₇0973"
$ws.Range("E50").Value = "  -7.50%  "
$ws.Range("D51").Value = "'0.0965"
$ws.Range("E51").Value = "  +2.57%  "
